# Update gh-pages to output generated at 456a3b4
# Apply updated 想去人数 (F column) counts across the relevant sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 308
$ws1.Range("F4").Value = 8126
$ws1.Range("F5").Value = 5922
$ws1.Range("F7").Value = 95
$ws1.Range("F8").Value = 16
$ws1.Range("F10").Value = 300
$ws1.Range("F11").Value = 512
$ws1.Range("F12").Value = 68

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 1

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 308
$ws4.Range("F4").Value = 8126
$ws4.Range("F5").Value = 5922
$ws4.Range("F7").Value = 95
$ws4.Range("F8").Value = 16
$ws4.Range("F10").Value = 300
$ws4.Range("F12").Value = 1
$ws4.Range("F15").Value = 512
$ws4.Range("F16").Value = 68
